# Updated symbol list (cryptos.xlsx) -- refreshed price/volume(1h) figures
# for rows 2-27 and 40-48, matching the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric- / percent-looking strings (e.g. "257.93",
    # "0.17%") are kept as literal text instead of being auto-parsed into a
    # number/percentage by Excel, then drop back to the default style so no
    # extra formatting is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "257.93"
Set-TextValue "E2" "0.17%"
Set-TextValue "D3" "27.01"
Set-TextValue "E3" "-0.28%"
Set-TextValue "D4" "4.668"
Set-TextValue "E4" "-4.78%"
Set-TextValue "D5" "0.05909"
Set-TextValue "E5" "-0.69%"
Set-TextValue "E6" "-0.68%"
Set-TextValue "D7" "0.8544"
Set-TextValue "E7" "-1.72%"
Set-TextValue "D8" "0.9485"
Set-TextValue "E8" "-1.57%"
Set-TextValue "D9" "0.1402"
Set-TextValue "E9" "-0.81%"
Set-TextValue "D10" "0.05079"
Set-TextValue "E10" "42.98%"
Set-TextValue "D11" "0.07099"
Set-TextValue "E11" "-1.11%"
Set-TextValue "D12" "0.03107"
Set-TextValue "D13" "0.09152"
Set-TextValue "E13" "-0.99%"
Set-TextValue "D14" "0.001524"
Set-TextValue "E14" "-1.57%"
Set-TextValue "D15" "0.0006030"
Set-TextValue "E15" "-0.51%"
Set-TextValue "D16" "0.006113"
Set-TextValue "E16" "2.31%"
Set-TextValue "E17" "0.45%"
Set-TextValue "E18" "-2.21%"
Set-TextValue "E19" "-0.04%"
Set-TextValue "D20" "0.3056"
Set-TextValue "E20" "-2.85%"
Set-TextValue "D21" "0.1277"
Set-TextValue "E21" "-2.27%"
Set-TextValue "D22" "3.826"
Set-TextValue "E22" "7.91%"
Set-TextValue "D23" "0.04264"
Set-TextValue "E23" "-0.37%"
Set-TextValue "E24" "-0.24%"
Set-TextValue "D25" "0.004296"
Set-TextValue "E25" "-4.90%"
Set-TextValue "E27" "29.88%"
Set-TextValue "D40" "0.03832"
Set-TextValue "E40" "0.00%"
Set-TextValue "D41" "0.006320"
Set-TextValue "E41" "58.36%"
Set-TextValue "E42" "-0.21%"
Set-TextValue "D43" "0.002200"
Set-TextValue "E43" "-5.20%"
Set-TextValue "D44" "0.01383"
Set-TextValue "E44" "31.39%"
Set-TextValue "D45" "0.00005402"
Set-TextValue "E45" "-1.61%"
Set-TextValue "D47" "0.05100"
Set-TextValue "D48" "0.2524"
Set-TextValue "E48" "11,673.99%"
